$d = $word.ActiveDocument

# Locate the end of the "LOB1018: Física I (Requisito fraco)" paragraph.
# Deletion should start right after its paragraph mark (i.e. at the
# following blank paragraph).
$r1 = $d.Content.Duplicate
$r1.Find.Execute("LOB1018: Física I (Requisito fraco)", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r1.MoveEnd(1, 1)   # wdCharacter: include the paragraph mark
$startDel = $r1.End

# Locate the end of the copyright paragraph ("© 2020 . Contact: ...").
# Deletion should end right after its paragraph mark.
$r2 = $d.Content.Duplicate
$r2.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
                  $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.MoveEnd(1, 1)   # wdCharacter: include the paragraph mark
$endDel = $r2.End

# Remove the blank paragraph, the "Ver no Jupiter..." paragraph and the
# copyright paragraph in one shot, leaving the following blank paragraph
# (and the page-break paragraph after it) untouched.
$delRange = $d.Range($startDel, $endDel)
$delRange.Delete()
